$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(298).Insert()

$ws.Range("A298").Value = 10
$ws.Range("B298").Value = "Vega Modelo de Temuco"
$ws.Range("C298").Value = "La Araucanía"
$ws.Range("D298").Value = 44746
$ws.Range("E298").Value = 9
$ws.Range("F298").Value = 100112017
$ws.Range("G298").Value = "Apio"
$ws.Range("H298").Value = "Americana (o)"
$ws.Range("I298").Value = "Primera"
$ws.Range("J298").Value = 140
$ws.Range("K298").Value = 8000
$ws.Range("L298").Value = 9000
$ws.Range("M298").Value = 8536
$ws.Range("N298").Value = "$/docena de matas"
$ws.Range("O298").Value = "Provincia del Elquí"
$ws.Range("P298").Value = 1423
$ws.Range("Q298").Value = 6
$ws.Range("R298").Value = "Hortaliza"
